$wb = $excel.ActiveWorkbook
$tasks = $wb.Worksheets.Item("Tasks")
$src = $tasks.Range("A6:M6")
$dst = $tasks.Range("A50:M50")
$src.Copy($dst)
Write-Host "copied row 6 to row 50"
